$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column holds values that look numeric (e.g. "43.218.30", "1.00", "0.0794")
# but must stay plain text, exactly as authored. Force Text format before writing
# so Excel does not silently reinterpret/round them as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Updated crypto price/volume data (rows 2-51: row, Coin, Link, Price, Volume(1h))
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '43.218.30', '  +1.95%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.312.43', '  +1.39%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  +0.09%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '301.40', '  +0.39%  '),
    @(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '99.32', '  +2.83%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.506', '  +0.91%  '),
    @(8, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  +0.03%  '),
    @(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.510', '  +2.39%  '),
    @(10, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '34.26', '  +3.76%  '),
    @(11, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.0794', '  +0.60%  '),
    @(12, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '49.17', '  -0.87%  '),
    @(13, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.116', '  +2.69%  '),
    @(14, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '17.47', '  +12.64%  '),
    @(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '6.81', '  +2.18%  '),
    @(16, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.673.98', '  +1.61%  '),
    @(17, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.355.03', '  +1.27%  '),
    @(18, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.816', '  +3.71%  '),
    @(19, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '43.171.40', '  +2.00%  '),
    @(20, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '11.89', '  +3.07%  '),
    @(21, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0₃0905', '  +0.74%  '),
    @(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.09', '  +1.12%  '),
    @(23, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '67.68', '  +1.53%  '),
    @(24, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '237.59', '  +1.62%  '),
    @(25, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '2.08', '  +8.06%  '),
    @(26, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  +0.07%  '),
    @(27, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '2.46', '  -0.67%  '),
    @(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '24.59', '  +1.72%  '),
    @(29, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.09', '  -7.31%  '),
    @(30, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '34.29', '  +0.91%  '),
    @(31, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '166.23', '  +0.90%  '),
    @(32, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.16', '  +0.31%  '),
    @(33, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '1.00', '  +0.12%  '),
    @(34, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.98', '  +0.11%  '),
    @(35, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '4.64', '  +5.69%  '),
    @(36, 'WEMIXToken', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '2.44', '  +3.42%  '),
    @(37, 'Celestia', 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia', '17.18', '  +6.26%  '),
    @(38, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0701', '  +0.75%  '),
    @(39, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.102', '  +3.81%  '),
    @(40, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.84', '  -0.07%  '),
    @(41, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.80', '  +3.05%  '),
    @(42, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.110', '  -0.58%  '),
    @(43, 'ApeXProtocol', 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex', '2.37', '  -3.20%  '),
    @(44, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.999.37', '  +2.02%  '),
    @(45, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0287', '  +1.55%  '),
    @(46, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '17.77', '  +0.58%  '),
    @(47, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '9.89', '  +1.65%  '),
    @(48, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '2.88', '  +2.27%  '),
    @(49, 'MultiversX', 'https://coinranking.com/coin/omwkOTglq+multiversx-egld', '53.88', '  +2.46%  '),
    @(50, 'RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '2.538.82', '  +1.38%  '),
    @(51, 'THORChain', 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune', '4.63', '  +0.35%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
